$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows before the existing row 134 (xpertacf_indigenous block),
# pushing the old rows 134-138 down to 139-143.
$ws.Rows("134:138").Insert()

# Populate the new rows' label cells in the order that reproduces the target
# shared-string table ordering (215..219):
#   215 econ_inflectioncost_xpertacf_indigenous
#   216 econ_startupcost_xpertacf_indigenous
#   217 econ_startupduration_xpertacf_indigenous
#   218 econ_saturation_xpertacf_indigenous
#   219 econ_unitcost_xpertacf_indigenous
$ws.Range("A135").Value = "econ_inflectioncost_xpertacf_indigenous"
$ws.Range("A136").Value = "econ_startupcost_xpertacf_indigenous"
$ws.Range("A137").Value = "econ_startupduration_xpertacf_indigenous"
$ws.Range("A138").Value = "econ_saturation_xpertacf_indigenous"
$ws.Range("A134").Value = "econ_unitcost_xpertacf_indigenous"

# Fill in the numeric values for the new rows.
$ws.Range("B134").Value = 30.26
$ws.Range("B135").Value = 0
$ws.Range("B136").Value = 662
$ws.Range("B137").Value = 1
$ws.Range("B138").Value = 0.9

# Update the active selection to match the saved view state.
$ws.Range("A135").Select() | Out-Null
